# Rename the existing sheet and add a second sheet after it, matching the
# target workbook layout: createCustomerData (sheetId 1) then
# openAccountData (sheetId 2), with openAccountData ending up the active tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "createCustomerData"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "openAccountData"

# --- createCustomerData: two more data rows ---------------------------
$ws1.Cells.Item(3, 1).Value = "John"
$ws1.Cells.Item(3, 2).Value = "wick"
$ws1.Cells.Item(4, 1).Value = "Malena"
$ws1.Cells.Item(4, 2).Value = "Anderson"
# leading apostrophe keeps these numeric-looking post codes stored as text
# (same quotePrefix style already used by the existing 623001 cell)
$ws1.Cells.Item(3, 3).Value = "'650001"
$ws1.Cells.Item(4, 3).Value = "'632102"
$ws1.Cells.Item(3, 4).Value = "Customer added successfully"
$ws1.Cells.Item(4, 4).Value = "Customer added successfully"

# --- openAccountData: new sheet content --------------------------------
$ws2.Cells.Item(1, 1).Value = "customerName"
$ws2.Cells.Item(1, 2).Value = "currency"
$ws2.Cells.Item(2, 2).Value = "Pound"
$ws2.Cells.Item(2, 1).Value = "Harry Potter"
$ws2.Cells.Item(1, 3).Value = "alert"
# quoted so the cell picks up the same quotePrefix style (s="1") the
# source workbook uses, matching the target markup
$ws2.Cells.Item(2, 3).Value = "'Account created successfully"

$ws2.Columns.Item(1).ColumnWidth = 13.86

# --- selections: leave createCustomerData parked on H22, and make
# openAccountData the active sheet selected at D7 -----------------------
$ws1.Range("H22").Select()
$ws2.Range("D7").Select()

Write-Host "done"
